$d = $word.ActiveDocument

# Replace "De en stelde" with "De klant stelde" in the paragraph about the
# odoo software link, effectively inserting the word "klant" that was
# missing from the sentence.
$d.Content.Find.Execute("De en stelde voor om eens te zoeken", $true, $false, $false, $false, $false,
                         $true, 1, $false, "De klant stelde voor om eens te zoeken", 2)
